$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 through 15 (these rows are removed entirely per the diff)
$ws.Range("A6:F15").EntireRow.Delete()

# Update row 2 values
$ws.Range("B2").Value = "NSE:RAMANEWS"
$ws.Range("C2").Value = "NSE:AMBUJACEM"
$ws.Range("D2").Value = "NSE:ESCORTS"
$ws.Range("E2").Value = "NSE:ADANIENSOL"
$ws.Range("F2").Value = "NSE:ANGELONE"

# Update row 3 values (B3 had content that must be cleared; D3 already blank, leave as-is)
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "NSE:JAGSNPHARM"
$ws.Range("E3").Value = "NSE:CROMPTON"
$ws.Range("F3").Value = "NSE:INDIGO"

# Update row 4 values (B4 had content that must be cleared; D4 already blank, leave as-is)
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "NSE:KIMS"
$ws.Range("E4").Value = "NSE:CUMMINSIND"
$ws.Range("F4").Value = "NSE:NMDC"

# Update row 5 values (B5, D5 already blank, leave as-is)
$ws.Range("C5").Value = "NSE:RALLIS"
$ws.Range("E5").Value = "NSE:DRREDDY"
$ws.Range("F5").Value = "NSE:OBEROIRLTY"
